$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("BTC/USDC", 103492.02516, 103550.78, 0.001, 0.05875483999999415, 0, 0.05875483999999415, 0.05677233575162764, "2025-05-17", "00:49:27", "2025-05-17", "01:00:00", 10.55609733333333),
    @("BTC/USDC", 103425.83464, 102826.53584, 0.001, -0.5992988000000041, 0, -0.5992988000000041, -0.5794478740113788, "2025-05-17", "02:26:39", "2025-05-17", "03:20:26", 53.77595986666667),
    @("BTC/USDC", 102974.83922, 103081.5932, 0.001, 0.1067539800000086, 0, 0.1067539800000086, 0.1036699652154199, "2025-05-17", "03:29:18", "2025-05-17", "04:25:04", 55.75855795)
)

$startRow = 14
# Entry Date (I) / Exit Date (K) look like dates and must be forced to text
# so Excel doesn't auto-convert them to date serials; time-only strings in
# Entry Time (J) / Exit Time (L) are not auto-converted, so no special
# handling is needed there.
$dateCols = @(9, 11)
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($r, $col)
        if ($dateCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $row[$c]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $row[$c]
        }
    }
}
